$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / Row 12 coin entries swapped positions (Bitrue <-> Mandala),
# plus refreshed price/volume figures across the sheet.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '304.40'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '-5.01%'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '39.78'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '-7.27%'
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '5.040'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '-2.26%'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '0.07667'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '-5.96%'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '4.254'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '-1.75%'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.582'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '-11.04%'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.8804'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '-7.49%'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.09767'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '-12.52%'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.1718'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '-7.16%'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.04433'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '-4.05%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.08887'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '-5.10%'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '-0.68%'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.001242'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '-3.39%'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.005844'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '-3.00%'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.354'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '-0.33%'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '2.422'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '-4.26%'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.035'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '-5.29%'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.1354'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '-2.54%'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '23.10%'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.04208'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '0.48%'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.001198'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '-4.46%'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.004051'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '-6.16%'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.0001224'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '9.90%'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '-0.21%'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.02316'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '-10.76%'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.05113'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '-7.61%'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.007977'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '1.81%'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '-5.43%'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.006517'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '-1.56%'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.001991'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '-6.11%'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.008490'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '0.37%'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.3020'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '-12.72%'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '-6.65%'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.00000000752'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '0.01%'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.007022'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '98.32%'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.003356'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '-3.38%'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.00002106'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '0.01%'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0002006'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '0.01%'
